$wb = $excel.ActiveWorkbook

foreach ($sheetName in @("展览", "全部类型")) {
    $ws = $wb.Worksheets.Item($sheetName)
    $ws.Range("F19").Value = 3751
    $ws.Range("F26").Value = 353
    $ws.Range("F28").Value = 1539
}
